$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52
$title52 = @'
Enriched line graph: A new structure for searching language collocations
'@
$pii52 = @'
S0960077920309012
'@
$abstract52 = @'
The specific terminology of a specialty language comes, essentially, from specific uses of already existing words and/or from specific combinations of words so called “collocations”. In this work we introduce a new mathematical structure (enriched line graph) and a new methodology to extract properties and characteristics of a type of multilayer linguistic networks associated with these types of languages. Specifically, this work is focused on the description of a methodology based on a variant of the PageRank algorithm to locate the linguistic collocations and on defining a new structure (enriched line graph) that can be interpreted as a certain type of “interpolation” between the original graph and its associated line graph, showing new results, properties and applications of this concept, and, in particular, certain characteristics of the specialty language produced by the scientific community of complex networks.
'@
$keywords52 = @'
Enriched line graph; Multilayer networks; PageRank; Interaction of words; Language collocations
'@
$ws.Range("B52").Value = $title52
$ws.Range("A52").Value = $pii52
$ws.Range("C52").Value = $abstract52
$ws.Range("D52").Value = $keywords52
$ws.Rows.Item(52).RowHeight = 165

# Row 53
$title53 = @'
Motifs and Motif Generalization in Chinese Word Networks
'@
$pii53 = @'
S1877050912001809
'@
$abstract53 = @'
The most signiﬁcant semantic unit of Chinese language is words composed of individual characters. This com–positional structure produces great variability and representability compared to individual characters, which is quite distinct from other languages. In this paper we utilized complex networks to model the composition of words from characters. We focus on network motifs, the local pattern which appears more often in a statistically signiﬁcant sense. Network motifs describe the most signiﬁcant connection pattern between these nodes. We investigated their functions and semantical relationship. We also investigated different generalizations of network motifs and analyzed the larger pattern in the complex networks. As the word network is quite huge and the motif detection is very slow when motifs are much larger, for larger pattern in the network we used topology generalization of simple motifs rather than carry out a thorough motif detection task. The results on motifs and motif generalization in this paper not only offer us a big picture how Chinese words are formed, but also support the conclusion that motifs play a very important role in research of complex systems.
'@
$keywords53 = @'
complex networks; motif; motif generalization
'@
$ws.Range("B53").Value = $title53
$ws.Range("A53").Value = $pii53
$ws.Range("C53").Value = $abstract53
$ws.Range("D53").Value = $keywords53
$ws.Rows.Item(53).RowHeight = 210

# Row 54
$title54 = @'
COVID-19 Symptoms app analysis to foresee healthcare impacts: Evidence from Northern Ireland
'@
$pii54 = @'
S1568494621011169
'@
$abstract54 = @'
Mobile health (mHealth) technologies, such as symptom tracking apps, are crucial for coping with the global pandemic crisis by providing near real-time, in situ information for the medical and governmental response. However, in such a dynamic and diverse environment, methods are still needed to support public health decision-making. This paper uses the lens of strong structuration theory to investigate networks of COVID-19 symptoms in the Belfast metropolitan area. A self-supervised machine learning method measuring information entropy was applied to the Northern Ireland COVIDCare app. The findings reveal: (1) relevant stratifications of disease symptoms, (2) particularities in health-wealth networks, and (3) the predictive potential of artificial intelligence to extract entangled knowledge from data in COVID-related apps. The proposed method proved to be effective for near real-time in-situ analysis of COVID-19 progression and to focus and complement public health decisions. Our contribution is relevant to an understanding of SARS-COV-2 symptom entanglements in localised environments. It can assist decision-makers in designing both reactive and proactive health measures that should be personalised to the heterogeneous needs of different populations. Moreover, near real-time assessment of pandemic symptoms using digital technologies will be critical to create early warning systems of emerging SARS-CoV-2 strains and predict the need for healthcare resources.
'@
$keywords54 = @'
COVID-19; SARS-COV-2; Strong structuration theory; Semantic networks; Mobile app; Location analytics; Symptoms assessment
'@
$ws.Range("B54").Value = $title54
$ws.Range("A54").Value = $pii54
$ws.Range("C54").Value = $abstract54
$ws.Range("D54").Value = $keywords54
$ws.Rows.Item(54).RowHeight = 255

# Row 55
$title55 = @'
A network intrusion detection method based on semantic Re-encoding and deep learning
'@
$pii55 = @'
S1084804520301624
'@
$abstract55 = @'
In recent years, with the increase of human activities in cyberspace, intrusion events, such as network penetration, detection and attack, tend to be frequent and hidden. The traditional intrusion detection methods which prefer rules are not enough to deal with the increasingly complex network intrusion flow. However, the generalization ability of intrusion detection system based on classical machine learning method is still insufficient, and the false alarm rate is high. Aiming at this problem, we consider that normal network traffic and intrusion network traffic are obviously different in several semantic dimensions, though the intrusion traffic is more and more covert. Then we propose a new intrusion detection method, named SRDLM, based on semantic re-encoding and deep learning. The SRDLM method re-encodes the semantics of network traffic, increases the distinguish ability of traffic, and enhances the generalization ability of the algorithm by using deep learning technology, thus effectively improving the accuracy and robustness of the algorithm. The accuracy of the SRDLC algorithm for Web character injection network attack detection is over 99%. When detecting the NSL-KDD data set, the average performance is improved by more than 8% compared with the traditional machine learning method.
'@
$keywords55 = @'
Intrusion detection; Semantic re-encoding; Deep learning
'@
$ws.Range("B55").Value = $title55
$ws.Range("A55").Value = $pii55
$ws.Range("C55").Value = $abstract55
$ws.Range("D55").Value = $keywords55
$ws.Rows.Item(55).RowHeight = 225

# Row 56
$title56 = @'
Metabolic networks classification and knowledge discovery by information granulation
'@
$pii56 = @'
S1476927119302440
'@
$abstract56 = @'
Graphs are powerful structures able to capture topological and semantic information from data, hence suitable for modelling a plethora of real-world (complex) systems. For this reason, graph-based pattern recognition gained a lot of attention in recent years. In this paper, a general-purpose classification system in the graphs domain is presented. When most of the information of the available patterns can be encoded in edge labels, an information granulation-based approach is highly discriminant and allows for the identification of semantically meaningful edges. The proposed classification system has been tested on the entire set of organisms (5299) for which metabolic networks are known, allowing for both a perfect mirroring of the underlying taxonomy and the identification of most discriminant metabolic reactions and pathways. The widespread diffusion of graph (network) structures in biology makes the proposed pattern recognition approach potentially very useful in many different fields of application. More specifically, the possibility to have a reliable metric to compare different metabolic systems is instrumental in emerging fields like microbiome analysis and, more in general, for proposing metabolic networks as a universal phenotype spanning the entire tree of life and in direct contact with environmental cues.
'@
$keywords56 = @'
Granular computing; Embedding spaces; Support vector machines; Computational biology; Metabolic pathways; Complex networks
'@
$ws.Range("B56").Value = $title56
$ws.Range("A56").Value = $pii56
$ws.Range("C56").Value = $abstract56
$ws.Range("D56").Value = $keywords56
$ws.Rows.Item(56).RowHeight = 240

# Row 57
$title57 = @'
Semantic networks based on titles of scientific papers
'@
$pii57 = @'
S0378437110010125
'@
$abstract57 = @'
In this paper we study the topological structure of semantic networks based on titles of papers published in scientific journals. It discusses its properties and presents some reflections on how the use of social and complex network models can contribute to the diffusion of knowledge. The proposed method presented here is applied to scientific journals where the titles of papers are in English or in Portuguese. We show that the topology of studied semantic networks are small-world and scale-free.
'@
$keywords57 = @'
Semantic networks; Complex networks; Social network analysis
'@
$ws.Range("B57").Value = $title57
$ws.Range("A57").Value = $pii57
$ws.Range("C57").Value = $abstract57
$ws.Range("D57").Value = $keywords57
$ws.Rows.Item(57).RowHeight = 90

# Row 58
$title58 = @'
Structure–semantics interplay in complex networks and its effects on the predictability of similarity in texts
'@
$pii58 = @'
S0378437112003044
'@
$abstract58 = @'
The classification of texts has become a major endeavor with so much electronic material available, for it is an essential task in several applications, including search engines and information retrieval. There are different ways to define similarity for grouping similar texts into clusters, as the concept of similarity may depend on the purpose of the task. For instance, in topic extraction similar texts mean those within the same semantic field, whereas in author recognition stylistic features should be considered. In this study, we introduce ways to classify texts employing concepts of complex networks, which may be able to capture syntactic, semantic and even pragmatic features. The interplay between various metrics of the complex networks is analyzed with three applications, namely identification of machine translation (MT) systems, evaluation of quality of machine translated texts and authorship recognition. We shall show that topological features of the networks representing texts can enhance the ability to identify MT systems in particular cases. For evaluating the quality of MT texts, on the other hand, high correlation was obtained with methods capable of capturing the semantics. This was expected because the golden standards used are themselves based on word co-occurrence. Notwithstanding, the Katz similarity, which involves semantic and structure in the comparison of texts, achieved the highest correlation with the NIST measurement, indicating that in some cases the combination of both approaches can improve the ability to quantify quality in MT. In authorship recognition, again the topological features were relevant in some contexts, though for the books and authors analyzed good results were obtained with semantic features as well. Because hybrid approaches encompassing semantic and topological features have not been extensively used, we believe that the methodology proposed here may be useful to enhance text classification considerably, as it combines well-established strategies.
'@
$keywords58 = @'
Similarity index; Complex networks; Machine translation evaluation; Topological analysis; Authorship recognition
'@
$ws.Range("B58").Value = $title58
$ws.Range("A58").Value = $pii58
$ws.Range("C58").Value = $abstract58
$ws.Range("D58").Value = $keywords58
$ws.Rows.Item(58).RowHeight = 360

# Row 59
$title59 = @'
Support for browsing in an intelligent text retrieval system
'@
$pii59 = @'
S0020737389800148
'@
$abstract59 = @'
Browsing is potentially an extremely important technique for retrieving text documents from large knowledge bases. The advantages of this technique are that users get immediate feedback from the structure of the knowledge base and exert complete control over the outcome of the search. The primary disadvantages are that it is easy to get lost in a complex network of nodes representing documents and concepts, and there is no guarantee that a browsing search will be as effective as a more conventional search. In this paper, we show how a browsing capability can be integrated into an intelligent text retrieval system. The disadvantages mentioned above are avoided by providing facilities for controlling the browsing and for using the information derived during browsing in more formal search strategies. The architecture of the text retrieval system is described and the browsing techniques are illustrated using an example session.
'@
$ws.Range("B59").Value = $title59
$ws.Range("A59").Value = $pii59
$ws.Range("C59").Value = $abstract59
$ws.Rows.Item(59).RowHeight = 45

# Row 60
$title60 = @'
RGB×D: Learning depth-weighted RGB patches for RGB-D indoor semantic segmentation
'@
$pii60 = @'
S0925231221011930
'@
$abstract60 = @'
Significant advances have been made in designing CNNs for RGB semantic segmentation. However, these CNNs are not widely adopted for RGB-D segmentation, due to the asymmetry between the RGB and depth modalities. Instead, dedicated architectures are designed to fuse them for effective RGB-D segmentation, wherein complex structures are often employed, resulting in much increased computational cost. In this paper, we propose a novel way to learn the fusion of RGB and depth information in an early stage. This enables our method to easily adopt existing RGB segmentation networks with minimal modification. Our method is simple yet effective to build a bridge between RGB and RGBD semantic segmentation, so as to avoid designing a far more complex network structure for RGBD segmentation. The proposed method treats RGB and depth information in an inherently asymmetric manner, and to the best of our knowledge, this is the first approach that learns to fuse them in a multiplicative manner for RGB-D segmentation; thus, we call it RGB×D. Extensive experiments and ablation studies on the challenging NYUDv2, SUN RGB-D and Cityscapes semantic segmentation benchmarks show that the proposed RGB×D offers a consistent improvement over several baselines.
'@
$keywords60 = @'
RGB-D indoor semantic segmentation; Depth information; Deep learning
'@
$ws.Range("B60").Value = $title60
$ws.Range("A60").Value = $pii60
$ws.Range("C60").Value = $abstract60
$ws.Range("D60").Value = $keywords60
$ws.Rows.Item(60).RowHeight = 225

# Row 61
$title61 = @'
Using virtual edges to improve the discriminability of co-occurrence text networks
'@
$pii61 = @'
S037843712030707X
'@
$abstract61 = @'
Word co-occurrence networks have been employed to analyze texts both in the practical and theoretical scenarios. Despite the relative success in several applications, traditional co-occurrence networks fail in establishing links between similar words whenever they appear distant in the text. Here we investigate whether the use of word embeddings as a tool to create virtual links in co-occurrence networks may improve the quality of classification systems. Our results revealed that the discriminability in the stylometry task is improved when using Glove, Word2Vec and FastText. In addition, we found that optimized results are obtained when stopwords are not disregarded and a simple global thresholding strategy is used to establish virtual links. Because the proposed approach is able to improve the representation of texts as complex networks, we believe that it could be extended to study other natural language processing tasks. Likewise, theoretical languages studies could benefit from the adopted enriched representation of word co-occurrence networks.
'@
$keywords61 = @'
Network science; Language networks; Text networks; Co-occurrence networks; Semantic networks; Word embeddings
'@
$ws.Range("B61").Value = $title61
$ws.Range("A61").Value = $pii61
$ws.Range("C61").Value = $abstract61
$ws.Range("D61").Value = $keywords61
$ws.Rows.Item(61).RowHeight = 180

$ws.Range("D61").Select()
